$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.578.75"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.577.81"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0892"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "1.803.18"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "1.561.85"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").Value = "28.581.12"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  -4.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("E29").Value = "  -2.52%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0485"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.77%  "
$ws.Range("E32").Value = "  -2.00%  "
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("D35").Value = "1.399.21"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +4.74%  "
$ws.Range("E37").Value = "  -4.11%  "
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("E39").Value = "  +2.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0166"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.522"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.961"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("D49").Value = "1.715.55"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("E51").Value = "  -1.20%  "
